# Sema3c-Plxnd1.xlsx : refresh with new TPM-derived NATMI output.
#
# The old table had 9 data rows (sending cluster = ECs/FAPs/MuSCs, each
# crossed with target cluster = ECs/FAPs/MuSCs). The refreshed run drops the
# "ECs" sending-cluster block entirely and recomputes the expression /
# specificity statistics for the remaining "FAPs" / "MuSCs" sending-cluster
# rows (6 rows total).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the obsolete "ECs" sending-cluster block (old rows 2:4). Deleting the
# range shifts the rows below it up, so the old "FAPs"/"MuSCs" blocks (old
# rows 5:10) become the new rows 2:7 - matching the refreshed table's shape
# (header + 6 data rows, dimension A1:T7).
$ws.Range("A2:T4").Delete()

# Re-assert the refreshed values for every surviving data row (columns A:T).
# Sending cluster / Ligand symbol / Receptor symbol / Target cluster
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Sema3c"
$ws.Range("C2").Value = "Plxnd1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 42.071953
$ws.Range("H2").Value = 126.215859
$ws.Range("I2").Value = 0.978774012990499
$ws.Range("J2").Value = 0.978774012990499
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 74.98881166666666
$ws.Range("N2").Value = 224.966435
$ws.Range("O2").Value = 0.6650357743745379
$ws.Range("P2").Value = 0.6650357743745379
$ws.Range("Q2").Value = 3154.925759965852
$ws.Range("R2").Value = 28394.33183969266
$ws.Range("S2").Value = 0.6509197336668106
$ws.Range("T2").Value = 0.6509197336668106

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Sema3c"
$ws.Range("C3").Value = "Plxnd1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 42.071953
$ws.Range("H3").Value = 126.215859
$ws.Range("I3").Value = 0.978774012990499
$ws.Range("J3").Value = 0.978774012990499
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.661646333333335
$ws.Range("N3").Value = 28.984939
$ws.Range("O3").Value = 0.08568398816056159
$ws.Range("P3").Value = 0.08568398816056158
$ws.Range("Q3").Value = 406.4843304386224
$ws.Range("R3").Value = 3658.358973947601
$ws.Range("S3").Value = 0.08386526094094328
$ws.Range("T3").Value = 0.08386526094094326

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Sema3c"
$ws.Range("C4").Value = "Plxnd1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 42.071953
$ws.Range("H4").Value = 126.215859
$ws.Range("I4").Value = 0.978774012990499
$ws.Range("J4").Value = 0.978774012990499
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 28.10860633333333
$ws.Range("N4").Value = 84.325819
$ws.Range("O4").Value = 0.2492802374649006
$ws.Range("P4").Value = 0.2492802374649006
$ws.Range("Q4").Value = 1182.583964551502
$ws.Range("R4").Value = 10643.25568096352
$ws.Range("S4").Value = 0.2439890183827453
$ws.Range("T4").Value = 0.2439890183827453

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Sema3c"
$ws.Range("C5").Value = "Plxnd1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.912385
$ws.Range("H5").Value = 2.737155
$ws.Range("I5").Value = 0.021225987009501
$ws.Range("J5").Value = 0.021225987009501
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 74.98881166666666
$ws.Range("N5").Value = 224.966435
$ws.Range("O5").Value = 0.6650357743745379
$ws.Range("P5").Value = 0.6650357743745379
$ws.Range("Q5").Value = 68.41866693249166
$ws.Range("R5").Value = 615.7680023924249
$ws.Range("S5").Value = 0.01411604070772738
$ws.Range("T5").Value = 0.01411604070772738

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Sema3c"
$ws.Range("C6").Value = "Plxnd1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.912385
$ws.Range("H6").Value = 2.737155
$ws.Range("I6").Value = 0.021225987009501
$ws.Range("J6").Value = 0.021225987009501
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.661646333333335
$ws.Range("N6").Value = 28.984939
$ws.Range("O6").Value = 0.08568398816056159
$ws.Range("P6").Value = 0.08568398816056158
$ws.Range("Q6").Value = 8.815141189838336
$ws.Range("R6").Value = 79.33627070854502
$ws.Range("S6").Value = 0.001818727219618317
$ws.Range("T6").Value = 0.001818727219618317

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Sema3c"
$ws.Range("C7").Value = "Plxnd1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.912385
$ws.Range("H7").Value = 2.737155
$ws.Range("I7").Value = 0.021225987009501
$ws.Range("J7").Value = 0.021225987009501
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 28.10860633333333
$ws.Range("N7").Value = 84.325819
$ws.Range("O7").Value = 0.2492802374649006
$ws.Range("P7").Value = 0.2492802374649006
$ws.Range("Q7").Value = 25.64587078943833
$ws.Range("R7").Value = 230.812837104945
$ws.Range("S7").Value = 0.005291219082155302
$ws.Range("T7").Value = 0.005291219082155302
